$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2197.755167519784
$ws.Cells.Item(2, 3).Value = 2206.848774439267
$ws.Cells.Item(2, 4).Value = 2214.805403650955
$ws.Cells.Item(2, 5).Value = 2212.565927724553
$ws.Cells.Item(3, 2).Value = 2202.387130431995
$ws.Cells.Item(3, 3).Value = 2209.291490143483
$ws.Cells.Item(3, 4).Value = 2215.442358444167
$ws.Cells.Item(3, 5).Value = 2213.057972484487
$ws.Cells.Item(4, 2).Value = 2180.219103371565
$ws.Cells.Item(4, 3).Value = 2194.109350309947
$ws.Cells.Item(4, 4).Value = 2206.561388682944
$ws.Cells.Item(4, 5).Value = 2205.988257890221
$ws.Cells.Item(5, 2).Value = 2199.509380308793
$ws.Cells.Item(5, 3).Value = 2206.671397389214
$ws.Cells.Item(5, 4).Value = 2213.921757419303
$ws.Cells.Item(5, 5).Value = 2212.657041173416
$ws.Cells.Item(6, 2).Value = 2207.103148948721
$ws.Cells.Item(6, 3).Value = 2213.597321349207
$ws.Cells.Item(6, 4).Value = 2219.610231407072
$ws.Cells.Item(6, 5).Value = 2217.150133649538
$ws.Cells.Item(7, 2).Value = 2192.770205129404
$ws.Cells.Item(7, 3).Value = 2201.945082437705
$ws.Cells.Item(7, 4).Value = 2210.892787299947
$ws.Cells.Item(7, 5).Value = 2210.270146677729
$ws.Cells.Item(8, 2).Value = 2189.482557421556
$ws.Cells.Item(8, 3).Value = 2197.242565447381
$ws.Cells.Item(8, 4).Value = 2204.46439533886
$ws.Cells.Item(8, 5).Value = 2202.406249820687
$ws.Cells.Item(9, 2).Value = 2195.799347806329
$ws.Cells.Item(9, 3).Value = 2202.093252104228
$ws.Cells.Item(9, 4).Value = 2208.723841640438
$ws.Cells.Item(9, 5).Value = 2207.082550820086
$ws.Cells.Item(10, 2).Value = 2056.801072003444
$ws.Cells.Item(10, 3).Value = 2111.90225829121
$ws.Cells.Item(10, 4).Value = 2161.726782239577
$ws.Cells.Item(10, 5).Value = 2178.023297273562
$ws.Cells.Item(11, 2).Value = 2044.274956698688
$ws.Cells.Item(11, 3).Value = 2111.254853176406
$ws.Cells.Item(11, 4).Value = 2160.859079728361
$ws.Cells.Item(11, 5).Value = 2177.081546657847
$ws.Cells.Item(12, 2).Value = 1880.46872151288
$ws.Cells.Item(12, 3).Value = 1993.751378557799
$ws.Cells.Item(12, 4).Value = 2098.669692755373
$ws.Cells.Item(12, 5).Value = 2135.293080965606
$ws.Cells.Item(13, 2).Value = 2043.082140312506
$ws.Cells.Item(13, 3).Value = 2094.715307719508
$ws.Cells.Item(13, 4).Value = 2150.349025874085
$ws.Cells.Item(13, 5).Value = 2170.737924912611
